# 2018.09.24 - Added Pipe Level to Level 1-1 (and imported all tiles)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calcs")

# The "Item Count" row for row 23 increases by one (Pipe Level tile added)
$ws.Range("A23").Value = 10

# New tile id 570 recorded for the Pipe Level column on row 23
$ws.Range("L23").Value = 570

# Reflect where the user was working when they made the edit
$ws.Range("A23").Select()
